$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed snapshot: swap ATOM/BTC rows (row4 becomes BTC, row5 becomes ATOM)
# and update the dollar/price figures for the latest pull of the dataframe.

# Row 3 - USDT
$ws.Range("D3").Value = 863.0471265
$ws.Range("E3").Value = 90.62048405
$ws.Range("F3").Value = 772.42664245
$ws.Range("H3").Value = 863.05

# Row 4 - now BTC
$ws.Range("B4").Value = "BTC"
$ws.Range("D4").Value = 0.01540148
$ws.Range("E4").Value = 0.01540148
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 33509
$ws.Range("H4").Value = 516.09

# Row 5 - now ATOM
$ws.Range("B5").Value = "ATOM"
$ws.Range("D5").Value = 24.0438
$ws.Range("E5").Value = 0.2372
$ws.Range("F5").Value = 23.8066
$ws.Range("G5").Value = 12.7
$ws.Range("H5").Value = 305.36

# Row 6 - ALGO
$ws.Range("G6").Value = 0.8508

# Row 7 - ETH
$ws.Range("G7").Value = 2211.73
